$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from Evaluation Statu" sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from Evaluation Statu")
$wsInclude.Name = "Include #0"

# 2. Update the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")

# 2a. URL value: pythia -> cicada
$wsMeta.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/eval-status"

# 2b. Date value updated
$wsMeta.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# 2c. Insert a new row 11 ("Jurisdiction") before the existing "Description" row,
#     copying the formatting of the row above it so the new row matches the
#     existing table styling.
$wsMeta.Rows.Item(11).Insert()
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsMeta.Range("A11").Value = "Jurisdiction"

# Use an existing empty-text cell as the copy source so B11 becomes a real
# (shared-string) empty-text cell rather than a truly blank cell.
$wsInclude.Range("A6").Copy()
$wsMeta.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# 3. Update the System URI value on the Include sheet: pythia -> cicada
$wsInclude.Range("B7").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/EvalStatus"
